# Target edit (per commit diff):
#   - Paragraph 1 "uji" becomes two runs "U" + "ji" (same paragraph).
#   - Two new paragraphs "Aku" and "Dia" are inserted after it.
#   - The trailing paragraph keeps only the _GoBack bookmark (no runs).
#
# Word inserts at paragraph 1 (index 0) were built from a COM-interop
# perspective: the paragraph mark / bookmark that sits at the very end of
# the story behaves specially with InsertParagraphAfter/Before, so the
# splits below are produced by inserting a literal carriage return
# (Chr 13, "`r") *before* that trailing position - this reliably keeps the
# bookmark attached to the paragraph mark while pushing it down to its own
# paragraph, exactly like a user placing the cursor there and pressing
# Enter.

$d = $word.ActiveDocument

# --- Step 1: split "uji" into "U" + "ji" as two separate runs -------------
# Temporarily break "u" into its own paragraph so editing its text does not
# get re-merged with "ji" (adjacent runs with identical formatting collapse
# back into one run). Then re-join the two paragraphs by deleting the
# paragraph mark between them, which leaves "U" and "ji" as distinct runs
# inside a single paragraph.
$first = $d.Range(0, 1)
$first.InsertParagraphAfter()

$uRange = $d.Range(0, 1)
$uRange.Text = "U"

$paraMark = $d.Range(1, 2)
$paraMark.Delete()

# --- Step 2: push the bookmark-only paragraph mark to its own paragraph ---
$tailStart = $d.Range(3, 3)
$tailStart.InsertBefore("`r")

# --- Step 3: insert the "Dia" paragraph right before the bookmark para ----
$diaBoundary = $d.Range(4, 4)
$diaBoundary.InsertBefore("Dia`r")

# --- Step 4: insert the "Aku" paragraph right before "Dia" ----------------
$akuBoundary = $d.Range(4, 4)
$akuBoundary.InsertBefore("Aku`r")

Write-Output ("Final text=[" + $d.Content.Text + "]")
Write-Output ("Paragraph count=" + $d.Paragraphs.Count)
